$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the H-column "position squared-error" formulas to take the
# square root, turning them into an actual (Euclidean) distance metric.
$ws.Range("H2").Formula = "=SQRT(((B2-E2)^2+(C2-F2)^2))"
$ws.Range("H3").Formula = "=SQRT(((B3-E3)^2+(C3-F3)^2))"
$ws.Range("H4:H12").Formula = "=SQRT(((B4-E4)^2+(C4-F4)^2))"

# Move the embedded scatter chart ("圖表 1") from its original spot
# (to the right of the data, starting around column I) down and over to
# start near column A beneath the data table.
$co = $ws.ChartObjects("圖表 1")
$co.Left = 48
$co.Top = 189
$co.Width = 939.216796875
$co.Height = 214.5

# Move the active selection to match where the author last left the
# cursor.
[void]$ws.Range("N14").Select()
